# "Generate Report for Handoff"
#
# The CI localization-status report is regenerated: the two tracked files
# (previously "5efaea19-...md" / "b351d048-...md") have been replaced by a
# fresh pair ("71ecf87f-...md" / "ffff454b1316-...md"), both of which are now
# freshly handed off (status "Ready for handoff") rather than handed back.
# Since a handoff has just happened, there is no "Latest Target File" /
# "Latest Handback File" yet, so those two columns are cleared out on the
# per-locale sheets.

$wb = $excel.ActiveWorkbook

$oldFile1 = "5efaea19-8b46-491c-b5f0-3fc48d58ff97"
$oldFile2 = "b351d048-6e0c-47da-b175-fa0e1a3d2857"
$newFile1 = "71ecf87f-1132-457e-9850-fc20d40adc5b"
$newFile2 = "ffff454b1316-8d5b-4ad0-816c-0da78165146d"
$newXlfHash = "debce0e461d240bb58e444abdca2ce992cfb69ba"

$status = "Ready for handoff"
$overviewDate = "2016-42-11 08:42:10"
$zhDate = "2016-03-11 08:42:07"
$deDate = "2016-03-11 08:42:10"
$emptyHandback = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = "$newFile1.md"
$ov.Range("B2").Value = $status
$ov.Range("C2").Value = $status
$ov.Range("D2").Value = $overviewDate

$ov.Range("A3").Value = "$newFile2.md"
$ov.Range("B3").Value = $status
$ov.Range("C3").Value = $status
$ov.Range("D3").Value = $overviewDate

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/$newFile1.md", "", "", "$newFile1.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/$newFile2.md", "", "", "$newFile2.md")

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Delete()

$zhXlf = "$newFile1.$newXlfHash.zh-cn.xlf"

$zh.Range("A2").Value = "$newFile1.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = $status
$zh.Range("D2").Value = $zhXlf
$zh.Range("E2").Value = $zhDate
$zh.Range("F2:G3").Clear()
$zh.Range("H2").Value = $emptyHandback
$zh.Range("I2").Value = "Include"

$zh.Range("A3").Value = "$newFile2.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = $status
$zh.Range("D3").Value = $zhXlf
$zh.Range("E3").Value = $zhDate
$zh.Range("H3").Value = $emptyHandback
$zh.Range("I3").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/$newFile1.md", "", "", "$newFile1.md")
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/$newFile1.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5faede3304f846602d1cca11ed3f74baec1c148/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", "", "", $zhXlf)

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/$newFile2.md", "", "", "$newFile2.md")
$zh.Hyperlinks.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/$newFile2.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5faede3304f846602d1cca11ed3f74baec1c148/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", "", "", $zhXlf)

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Delete()

$deXlf = "$newFile1.$newXlfHash.de-de.xlf"

$de.Range("A2").Value = "$newFile1.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = $status
$de.Range("D2").Value = $deXlf
$de.Range("E2").Value = $deDate
$de.Range("F2:G3").Clear()
$de.Range("H2").Value = $emptyHandback
$de.Range("I2").Value = "Include"

$de.Range("A3").Value = "$newFile2.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = $status
$de.Range("D3").Value = $deXlf
$de.Range("E3").Value = $deDate
$de.Range("H3").Value = $emptyHandback
$de.Range("I3").Value = "Include"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/$newFile1.md", "", "", "$newFile1.md")
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/$newFile1.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/69340c41d8f5df6026ddcbbe2603c08d4516ede0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", "", "", $deXlf)

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/$newFile2.md", "", "", "$newFile2.md")
$de.Hyperlinks.Add($de.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/$newFile2.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/69340c41d8f5df6026ddcbbe2603c08d4516ede0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", "", "", $deXlf)
